$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 291.1111
$ws.Range("I33").Value = 322.5
$ws.Range("K33").Value = 322.5
$ws.Range("M33").Value = -93.5
$ws.Range("H42").Value = 270.4
$ws.Range("I42").Value = 147.5
$ws.Range("J42").Value = 352.33334
$ws.Range("K42").Value = 442.5
$ws.Range("L42").Value = 1057.00002
$ws.Range("M42").Value = -212.5
$ws.Range("N42").Value = -1517.00002
$ws.Range("H62").Value = 3755.7856
$ws.Range("I62").Value = 3126.6365
$ws.Range("J62").Value = 6062.6665
$ws.Range("K62").Value = 3126.6365
$ws.Range("L62").Value = 6062.6665
$ws.Range("M62").Value = -2502.6365
$ws.Range("N62").Value = -7310.6665
$ws.Range("H65").Value = 3755.7856
$ws.Range("I65").Value = 3126.6365
$ws.Range("J65").Value = 6062.6665
$ws.Range("K65").Value = 15633.1825
$ws.Range("L65").Value = 30313.3325
$ws.Range("M65").Value = -12513.1825
$ws.Range("N65").Value = -36553.3325
$ws.Range("H69").Value = 40000
$ws.Range("J69").Value = 40000
$ws.Range("L69").Value = 120000
$ws.Range("N69").Value = -121748
$ws.Range("H72").Value = 40000
$ws.Range("J72").Value = 40000
$ws.Range("L72").Value = 360000
$ws.Range("N72").Value = -368736
$ws.Range("H98").Value = 2167.9678
$ws.Range("I98").Value = 2235.6072
$ws.Range("K98").Value = 2235.6072
$ws.Range("M98").Value = -737.6071999999999
$ws.Range("H122").Value = 2167.9678
$ws.Range("I122").Value = 2235.6072
$ws.Range("K122").Value = 6706.821599999999
$ws.Range("M122").Value = -4256.821599999999
$ws.Range("H132").Value = 7213.25
$ws.Range("I132").Value = 7540.316
$ws.Range("J132").Value = 999
$ws.Range("K132").Value = 22620.948
$ws.Range("L132").Value = 2997
$ws.Range("M132").Value = -20090.948
$ws.Range("N132").Value = -8057
$ws.Range("H137").Value = 2646.5217
$ws.Range("I137").Value = 1706.0714
$ws.Range("K137").Value = 5118.2142
$ws.Range("M137").Value = -2568.2142

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1366855.5
$ws.Range("I32").Value = 640067.5600000001
$ws.Range("J32").Value = 11905280
$ws.Range("K32").Value = 640067.5600000001
$ws.Range("L32").Value = 11905280
$ws.Range("M32").Value = -639780.5600000001
$ws.Range("N32").Value = -11905854
$ws.Range("H61").Value = 3474.3157
$ws.Range("I61").Value = 3335
$ws.Range("J61").Value = 3599.7
$ws.Range("K61").Value = 3335
$ws.Range("L61").Value = 3599.7
$ws.Range("M61").Value = -3123
$ws.Range("N61").Value = -4023.7
$ws.Range("H122").Value = 3372.923
$ws.Range("I122").Value = 3035.4285
$ws.Range("K122").Value = 9106.2855
$ws.Range("M122").Value = -6656.2855
$ws.Range("H136").Value = 3474.3157
$ws.Range("I136").Value = 3335
$ws.Range("J136").Value = 3599.7
$ws.Range("K136").Value = 10005
$ws.Range("L136").Value = 10799.1
$ws.Range("M136").Value = -7455
$ws.Range("N136").Value = -15899.1

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H106").Value = 27617
$ws.Range("J106").Value = 27617
$ws.Range("L106").Value = 27617
$ws.Range("N106").Value = -30141

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8340388
$ws.Range("I31").Value = 2611.6667
$ws.Range("J31").Value = 20847052
$ws.Range("K31").Value = 2611.6667
$ws.Range("L31").Value = 20847052
$ws.Range("M31").Value = -2316.6667
$ws.Range("N31").Value = -20847642
$ws.Range("H34").Value = 8340388
$ws.Range("I34").Value = 2611.6667
$ws.Range("J34").Value = 20847052
$ws.Range("K34").Value = 2611.6667
$ws.Range("L34").Value = 20847052
$ws.Range("M34").Value = -2409.6667
$ws.Range("N34").Value = -20847456
$ws.Range("H86").Value = 6489.4546
$ws.Range("I86").Value = 5597.4443
$ws.Range("J86").Value = 10503.5
$ws.Range("K86").Value = 5597.4443
$ws.Range("L86").Value = 10503.5
$ws.Range("M86").Value = -4474.4443
$ws.Range("N86").Value = -12749.5
$ws.Range("H89").Value = 6489.4546
$ws.Range("I89").Value = 5597.4443
$ws.Range("J89").Value = 10503.5
$ws.Range("K89").Value = 27987.2215
$ws.Range("L89").Value = 52517.5
$ws.Range("M89").Value = -22371.2215
$ws.Range("N89").Value = -63749.5
$ws.Range("H99").Value = 2741.4167
$ws.Range("I99").Value = 1848.5454
$ws.Range("J99").Value = 3496.923
$ws.Range("K99").Value = 1848.5454
$ws.Range("L99").Value = 3496.923
$ws.Range("M99").Value = -350.5454
$ws.Range("N99").Value = -6492.923
$ws.Range("H105").Value = 1819.05
$ws.Range("I105").Value = 1846.5294
$ws.Range("J105").Value = 1663.3334
$ws.Range("K105").Value = 1846.5294
$ws.Range("L105").Value = 1663.3334
$ws.Range("M105").Value = -99.5293999999999
$ws.Range("N105").Value = -5157.3334
$ws.Range("H126").Value = 2741.4167
$ws.Range("I126").Value = 1848.5454
$ws.Range("J126").Value = 3496.923
$ws.Range("K126").Value = 5545.6362
$ws.Range("L126").Value = 10490.769
$ws.Range("M126").Value = -3075.6362
$ws.Range("N126").Value = -15430.769
$ws.Range("H132").Value = 4432.758
$ws.Range("I132").Value = 3276.3076
$ws.Range("K132").Value = 9828.9228
$ws.Range("M132").Value = -7298.9228
$ws.Range("H134").Value = 3352.0571
$ws.Range("I134").Value = 3268.5483
$ws.Range("K134").Value = 9805.644899999999
$ws.Range("M134").Value = -7270.644899999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 3280.4443
$ws.Range("I140").Value = 3253
$ws.Range("K140").Value = 9759
$ws.Range("M140").Value = -4579

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5511.0625
$ws.Range("I70").Value = 4466.8887
$ws.Range("K70").Value = 4466.8887
$ws.Range("M70").Value = -4196.8887
$ws.Range("H73").Value = 5511.0625
$ws.Range("I73").Value = 4466.8887
$ws.Range("K73").Value = 4466.8887
$ws.Range("M73").Value = -3530.8887
$ws.Range("H122").Value = 3849878.8
$ws.Range("I122").Value = 6995959.5
$ws.Range("K122").Value = 20987878.5
$ws.Range("M122").Value = -20985428.5
$ws.Range("H132").Value = 2252.3872
$ws.Range("I132").Value = 2156.45
$ws.Range("K132").Value = 6469.349999999999
$ws.Range("M132").Value = -3939.349999999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 787.0909
$ws.Range("I16").Value = 787.0909
$ws.Range("K16").Value = 787.0909
$ws.Range("M16").Value = -617.0909
$ws.Range("H22").Value = 1648.9
$ws.Range("I22").Value = 1213
$ws.Range("J22").Value = 2666
$ws.Range("K22").Value = 1213
$ws.Range("L22").Value = 2666
$ws.Range("M22").Value = -918
$ws.Range("N22").Value = -3256
$ws.Range("H27").Value = 1648.9
$ws.Range("I27").Value = 1213
$ws.Range("J27").Value = 2666
$ws.Range("K27").Value = 1213
$ws.Range("L27").Value = 2666
$ws.Range("M27").Value = -1106
$ws.Range("N27").Value = -2880
$ws.Range("H132").Value = 16999
$ws.Range("I132").Value = 16999
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 50997
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -48467
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 4279.5835
$ws.Range("I136").Value = 3746.6667
$ws.Range("K136").Value = 11240.0001
$ws.Range("M136").Value = -8690.000100000001
$ws.Range("H140").Value = 102747.875
$ws.Range("J140").Value = 102747.875
$ws.Range("L140").Value = 102747.875
$ws.Range("N140").Value = -113107.875
$ws.Range("H141").Value = 99332.164
$ws.Range("J141").Value = 99332.164
$ws.Range("L141").Value = 99332.164
$ws.Range("N141").Value = -109692.164

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 6723.1763
$ws.Range("I81").Value = 6215.3335
$ws.Range("J81").Value = 7000.1816
$ws.Range("K81").Value = 12430.667
$ws.Range("L81").Value = 14000.3632
$ws.Range("M81").Value = -11369.667
$ws.Range("N81").Value = -16122.3632
$ws.Range("H84").Value = 6723.1763
$ws.Range("I84").Value = 6215.3335
$ws.Range("J84").Value = 7000.1816
$ws.Range("K84").Value = 62153.335
$ws.Range("L84").Value = 70001.81599999999
$ws.Range("M84").Value = -56849.335
$ws.Range("N84").Value = -80609.81599999999
$ws.Range("H136").Value = 3168.647
$ws.Range("I136").Value = 3419.2144
$ws.Range("K136").Value = 10257.6432
$ws.Range("M136").Value = -7707.643199999999
